# Disaggregation of commodity Copper
# 1) Rename the commodity label from "Copper ores and concentrates" to "Copper"
#    (the label sits in cell C4 of every year sheet, backed by a single shared
#    string, so update it everywhere it appears).
# 2) A handful of year sheets get a last-significant-digit recalculation of the
#    value in D4 that goes along with the relabeling/disaggregation.

$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    if ($ws.Cells.Item(4, 3).Value2 -eq "Copper ores and concentrates") {
        $ws.Cells.Item(4, 3).Value = "Copper"
    }
}

$d4Updates = @(
    @{ Year = [string]"2033"; Value = 110830.1039065614 },
    @{ Year = [string]"2039"; Value = 216811.3829355027 },
    @{ Year = [string]"2041"; Value = 278380.1093116245 },
    @{ Year = [string]"2045"; Value = 642552.158481146 },
    @{ Year = [string]"2067"; Value = 748329.7765664503 },
    @{ Year = [string]"2069"; Value = 939284.8480597934 },
    @{ Year = [string]"2072"; Value = 1396816.716286596 }
)

foreach ($entry in $d4Updates) {
    $ws = $wb.Worksheets.Item($entry.Year)
    $ws.Cells.Item(4, 4).Value = $entry.Value
}
